# Group2_Milestone2_PartA.xlsx - "Content inventory change and persona added"
#
# 1) Add two hyperlinked asset-source URLs (cells D10, D11, D13) that link
#    to the team's GitHub asset folders (images/originals and docs).
# 2) Remove the now-unused "Quantity" column (old column F) from the
#    inventory table, shifting Status/Comments left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Hyperlinks: cell text becomes the URL itself, Excel auto-applies
#        the built-in "Link" cell style (underline, theme hyperlink color).
$ws.Hyperlinks.Add($ws.Range("D10"), "https://github.com/HIT226/Code-Fair-Assets/tree/master/images/originals")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://github.com/HIT226/Code-Fair-Assets/tree/master/docs")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://github.com/HIT226/Code-Fair-Assets/tree/master/docs")

# --- 2) Drop the "Quantity" column entirely (it was column F).
$ws.Columns("F").Delete()
